$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 4
$ws.Range("J4").Value = 4.45
$ws.Range("K4").Value = 2.07
$ws.Range("N4").Value = 6.6
$ws.Range("U4").Value = 1.88
$ws.Range("V4").Value = 1.82
$ws.Range("X4").Value = 22
$ws.Range("AC4").Value = 6.6
$ws.Range("AP4").Value = 29
$ws.Range("BA4").Value = 70
$ws.Range("BB4").Value = 250
